{"js": "// Applies the diff: updates the date line and the 25 multiplication\n// problems/answers inside the single table, matching the exact\n// old-text -> new-text pairs from the commit's OOXML diff.\nconst replacements = [\n  [\"2025-09-08 Monday\", \"2025-09-09 Tuesday\"],\n  [\"571\u00d79=5139\", \"548\u00d73=1644\"],\n  [\"216\u00d75=1080\", \"439\u00d74=1756\"],\n  [\"863\u00d78=6904\", \"987\u00d77=6909\"],\n  [\"431\u00d74=1724\", \"483\u00d74=1932\"],\n  [\"838\u00d79=7542\", \"486\u00d72=972\"],\n  [\"260\u00d73=780\", \"878\u00d78=7024\"],\n  [\"822\u00d77=5754\", \"983\u00d72=1966\"],\n  [\"523\u00d73=1569\", \"281\u00d75=1405\"],\n  [\"180\u00d78=1440\", \"671\u00d73=2013\"],\n  [\"827\u00d78=6616\", \"282\u00d72=564\"],\n  [\"123\u00d78=984\", \"852\u00d74=3408\"],\n  [\"414\u00d75=2070\", \"952\u00d72=1904\"],\n  [\"601\u00d75=3005\", \"597\u00d74=2388\"],\n  [\"642\u00d76=3852\", \"977\u00d75=4885\"],\n  [\"874\u00d75=4370\", \"977\u00d73=2931\"],\n  [\"282\u00d75=1410\", \"280\u00d75=1400\"],\n  [\"349\u00d79=3141\", \"709\u00d76=4254\"],\n  [\"975\u00d76=5850\", \"549\u00d78=4392\"],\n  [\"224\u00d75=1120\", \"932\u00d79=8388\"],\n  [\"213\u00d72=426\", \"685\u00d75=3425\"],\n  [\"834\u00d73=2502\", \"320\u00d79=2880\"],\n  [\"877\u00d77=6139\", \"907\u00d77=6349\"],\n  [\"496\u00d74=1984\", \"637\u00d72=1274\"],\n  [\"721\u00d76=4326\", \"477\u00d73=1431\"],\n  [\"891\u00d79=8019\", \"995\u00d73=2985\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found, cannot apply replacement: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Applies the diff: updates the date line and the 25 multiplication\n# problems/answers inside the single table, matching the exact\n# old-text -> new-text pairs from the commit's OOXML diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-09-08 Monday\", \"2025-09-09 Tuesday\"),\n    @(\"571\u00d79=5139\", \"548\u00d73=1644\"),\n    @(\"216\u00d75=1080\", \"439\u00d74=1756\"),\n    @(\"863\u00d78=6904\", \"987\u00d77=6909\"),\n    @(\"431\u00d74=1724\", \"483\u00d74=1932\"),\n    @(\"838\u00d79=7542\", \"486\u00d72=972\"),\n    @(\"260\u00d73=780\", \"878\u00d78=7024\"),\n    @(\"822\u00d77=5754\", \"983\u00d72=1966\"),\n    @(\"523\u00d73=1569\", \"281\u00d75=1405\"),\n    @(\"180\u00d78=1440\", \"671\u00d73=2013\"),\n    @(\"827\u00d78=6616\", \"282\u00d72=564\"),\n    @(\"123\u00d78=984\", \"852\u00d74=3408\"),\n    @(\"414\u00d75=2070\", \"952\u00d72=1904\"),\n    @(\"601\u00d75=3005\", \"597\u00d74=2388\"),\n    @(\"642\u00d76=3852\", \"977\u00d75=4885\"),\n    @(\"874\u00d75=4370\", \"977\u00d73=2931\"),\n    @(\"282\u00d75=1410\", \"280\u00d75=1400\"),\n    @(\"349\u00d79=3141\", \"709\u00d76=4254\"),\n    @(\"975\u00d76=5850\", \"549\u00d78=4392\"),\n    @(\"224\u00d75=1120\", \"932\u00d79=8388\"),\n    @(\"213\u00d72=426\", \"685\u00d75=3425\"),\n    @(\"834\u00d73=2502\", \"320\u00d79=2880\"),\n    @(\"877\u00d77=6139\", \"907\u00d77=6349\"),\n    @(\"496\u00d74=1984\", \"637\u00d72=1274\"),\n    @(\"721\u00d76=4326\", \"477\u00d73=1431\"),\n    @(\"891\u00d79=8019\", \"995\u00d73=2985\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found, cannot apply replacement: $oldText\"\n    }\n}\n"}
